$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 112266830
$ws.Range("B6").Value = 102158
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "CR"
$ws.Range("E6").Value = 1563
$ws.Range("F6").Value = "Bohuslind"
$ws.Range("G6").Value = "Tilia platyphyllos"
$ws.Range("H6").Value = "Scop."
$ws.Range("P6").Value = "Skäggstorp, Boh"
$ws.Range("Q6").Value = 313436
$ws.Range("R6").Value = 6428463
$ws.Range("S6").Value = 50
$ws.Range("T6").Value = "Västra Götaland"
$ws.Range("U6").Value = "Kungälv"
$ws.Range("V6").Value = "Bohuslän"
$ws.Range("W6").Value = "Solberga"

# Force these two as literal text (not auto-converted to a date serial)
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-08-31"
$ws.Range("Y6").NumberFormat = "General"
$ws.Range("Y6").Style = "Normal"

$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-08-31"
$ws.Range("AA6").NumberFormat = "General"
$ws.Range("AA6").Style = "Normal"

$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = "Bengt Westman"
$ws.Range("AX6").Value = "Bengt Westman"
